# "Generate Report for Handback" — update the localization-status report:
#  - flip the per-file status from "Ready for handoff" to
#    "Handed back: in sync with en-US" on every sheet that shows it
#  - add "Latest Target File" (F) / "Latest Handback File" (G) hyperlink
#    cells for each data row on the zh-cn and de-de sheets
#  - stamp the handback datetime in column H (per-language sheet)

$wb = $excel.ActiveWorkbook

$mdName  = "5995f510-1401-4e86-a5b1-e263820d3f24.md"
$zhXlf   = "5995f510-1401-4e86-a5b1-e263820d3f24.0badd6d7a8fddbde74b09176fe144e2cd2da3917.zh-cn.xlf"
$deXlf   = "5995f510-1401-4e86-a5b1-e263820d3f24.0badd6d7a8fddbde74b09176fe144e2cd2da3917.de-de.xlf"

$mdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/6233732112122934567495c5c7fe7e373d7e68ae/e2e/5995f510-1401-4e86-a5b1-e263820d3f24.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f6053283774a70ecf78488a5e7a0a29e6d8efb7b/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/5995f510-1401-4e86-a5b1-e263820d3f24.0badd6d7a8fddbde74b09176fe144e2cd2da3917.zh-cn.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/afafa7a8353d32a967c03a8f59c59cd18de4b889/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/5995f510-1401-4e86-a5b1-e263820d3f24.0badd6d7a8fddbde74b09176fe144e2cd2da3917.de-de.xlf"

$statusText = "Handed back: in sync with en-US"
$hyperlinkColor = 15570276   # BGR for FF6495ED, matching the sheet's existing hyperlink font

function Set-HandbackLink($ws, $cellRef, $text, $url) {
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, "", "", $text) | Out-Null
}

# --- Overview sheet: status column for both languages, both rows ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

Set-HandbackLink $wsZh "F2" $mdName $mdUrl
Set-HandbackLink $wsZh "G2" $zhXlf $zhXlfUrl
Set-HandbackLink $wsZh "F3" $mdName $mdUrl
Set-HandbackLink $wsZh "G3" $zhXlf $zhXlfUrl

foreach ($ref in @("F2", "G2", "F3", "G3")) {
    $wsZh.Range($ref).Font.Color = $hyperlinkColor
    $wsZh.Range($ref).Font.Underline = 2
}

$wsZh.Range("H2").Value = "2016-03-25 08:44:25"
$wsZh.Range("H3").Value = "2016-03-25 08:44:25"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

Set-HandbackLink $wsDe "F2" $mdName $mdUrl
Set-HandbackLink $wsDe "G2" $deXlf $deXlfUrl
Set-HandbackLink $wsDe "F3" $mdName $mdUrl
Set-HandbackLink $wsDe "G3" $deXlf $deXlfUrl

foreach ($ref in @("F2", "G2", "F3", "G3")) {
    $wsDe.Range($ref).Font.Color = $hyperlinkColor
    $wsDe.Range($ref).Font.Underline = 2
}

$wsDe.Range("H2").Value = "2016-03-25 08:44:40"
$wsDe.Range("H3").Value = "2016-03-25 08:44:40"
